$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- 1. Clear the manual line-spacing reduction on the autofit body
#        (<a:normAutofit lnSpcReduction="10000"/> -> <a:normAutofit/>)
$tf.AutoSize = 2

# --- 2. "Earth & Terra Satellite:" paragraph -> prepend "  Earth " as new runs,
#        leaving "& Terra Satellite:" as the trailing (original) run.
$para1 = $tr.Paragraphs(1, 1)
$ins1 = $para1.InsertBefore("  ")
$spacesRange = $tr.Characters($ins1.Start, 2)
$spacesRange.LanguageID = "pt-BR"

$para1Now = $tr.Paragraphs(1, 1)
$earthRange = $tr.Characters($para1Now.Start + 2, 6)
$earthRange.Text = "Earth "

# --- 3. Delete the sciencedirect reference paragraph entirely.
$tr.Paragraphs(9, 1).Delete()

# --- 4. Split the arxiv reference paragraph into a "- " run and a bare-URL run.
$arxivPara = $tr.Paragraphs(9, 1)
$prefixRange = $tr.Characters($arxivPara.Start, 2)
$prefixRange.Text = ""
$arxivParaNow = $tr.Paragraphs(9, 1)
$arxivParaNow.InsertBefore("- ")
